$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ntf3"
$ws.Cells.Item(2, 3).Value = "Ntrk3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 8.859944
$ws.Cells.Item(2, 8).Value = 26.579832
$ws.Cells.Item(2, 9).Value = 0.4921128329655918
$ws.Cells.Item(2, 10).Value = 0.4921128329655918
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.009008
$ws.Cells.Item(2, 14).Value = 0.027024
$ws.Cells.Item(2, 15).Value = 0.009775433435787767
$ws.Cells.Item(2, 16).Value = 0.009775433435787767
$ws.Cells.Item(2, 17).Value = 0.079810375552
$ws.Cells.Item(2, 18).Value = 0.718293379968
$ws.Cells.Item(2, 19).Value = 0.004810616241552087
$ws.Cells.Item(2, 20).Value = 0.004810616241552087

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ntf3"
$ws.Cells.Item(3, 3).Value = "Ntrk3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 8.859944
$ws.Cells.Item(3, 8).Value = 26.579832
$ws.Cells.Item(3, 9).Value = 0.4921128329655918
$ws.Cells.Item(3, 10).Value = 0.4921128329655918
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.07403633333333333
$ws.Cells.Item(3, 14).Value = 0.222109
$ws.Cells.Item(3, 15).Value = 0.08034383307391152
$ws.Cells.Item(3, 16).Value = 0.08034383307391152
$ws.Cells.Item(3, 17).Value = 0.6559577672986666
$ws.Cells.Item(3, 18).Value = 5.903619905688
$ws.Cells.Item(3, 19).Value = 0.03953823130531721
$ws.Cells.Item(3, 20).Value = 0.03953823130531721

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ntf3"
$ws.Cells.Item(4, 3).Value = "Ntrk3"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 8.859944
$ws.Cells.Item(4, 8).Value = 26.579832
$ws.Cells.Item(4, 9).Value = 0.4921128329655918
$ws.Cells.Item(4, 10).Value = 0.4921128329655918
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.8384493333333333
$ws.Cells.Item(4, 14).Value = 2.515348
$ws.Cells.Item(4, 15).Value = 0.9098807334903006
$ws.Cells.Item(4, 16).Value = 0.9098807334903007
$ws.Cells.Item(4, 17).Value = 7.428614140170667
$ws.Cells.Item(4, 18).Value = 66.85752726153599
$ws.Cells.Item(4, 19).Value = 0.4477639854187224
$ws.Cells.Item(4, 20).Value = 0.4477639854187225

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ntf3"
$ws.Cells.Item(5, 3).Value = "Ntrk3"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 5.867977666666667
$ws.Cells.Item(5, 8).Value = 17.603933
$ws.Cells.Item(5, 9).Value = 0.3259283708025871
$ws.Cells.Item(5, 10).Value = 0.3259283708025871
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.009008
$ws.Cells.Item(5, 14).Value = 0.027024
$ws.Cells.Item(5, 15).Value = 0.009775433435787767
$ws.Cells.Item(5, 16).Value = 0.009775433435787767
$ws.Cells.Item(5, 17).Value = 0.05285874282133334
$ws.Cells.Item(5, 18).Value = 0.475728685392
$ws.Cells.Item(5, 19).Value = 0.003186091093615443
$ws.Cells.Item(5, 20).Value = 0.003186091093615443

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Ntf3"
$ws.Cells.Item(6, 3).Value = "Ntrk3"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 5.867977666666667
$ws.Cells.Item(6, 8).Value = 17.603933
$ws.Cells.Item(6, 9).Value = 0.3259283708025871
$ws.Cells.Item(6, 10).Value = 0.3259283708025871
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.07403633333333333
$ws.Cells.Item(6, 14).Value = 0.222109
$ws.Cells.Item(6, 15).Value = 0.08034383307391152
$ws.Cells.Item(6, 16).Value = 0.08034383307391152
$ws.Cells.Item(6, 17).Value = 0.4344435505218889
$ws.Cells.Item(6, 18).Value = 3.909991954697
$ws.Cells.Item(6, 19).Value = 0.02618633461781499
$ws.Cells.Item(6, 20).Value = 0.02618633461781499

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Ntf3"
$ws.Cells.Item(7, 3).Value = "Ntrk3"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 5.867977666666667
$ws.Cells.Item(7, 8).Value = 17.603933
$ws.Cells.Item(7, 9).Value = 0.3259283708025871
$ws.Cells.Item(7, 10).Value = 0.3259283708025871
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.8384493333333333
$ws.Cells.Item(7, 14).Value = 2.515348
$ws.Cells.Item(7, 15).Value = 0.9098807334903006
$ws.Cells.Item(7, 16).Value = 0.9098807334903007
$ws.Cells.Item(7, 17).Value = 4.920001962631556
$ws.Cells.Item(7, 18).Value = 44.280017663684
$ws.Cells.Item(7, 19).Value = 0.2965559450911566
$ws.Cells.Item(7, 20).Value = 0.2965559450911566

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Ntf3"
$ws.Cells.Item(8, 3).Value = "Ntrk3"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 3.275965666666667
$ws.Cells.Item(8, 8).Value = 9.827897
$ws.Cells.Item(8, 9).Value = 0.1819587962318212
$ws.Cells.Item(8, 10).Value = 0.1819587962318212
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.009008
$ws.Cells.Item(8, 14).Value = 0.027024
$ws.Cells.Item(8, 15).Value = 0.009775433435787767
$ws.Cells.Item(8, 16).Value = 0.009775433435787767
$ws.Cells.Item(8, 17).Value = 0.02950989872533333
$ws.Cells.Item(8, 18).Value = 0.265589088528
$ws.Cells.Item(8, 19).Value = 0.001778726100620238
$ws.Cells.Item(8, 20).Value = 0.001778726100620238

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Ntf3"
$ws.Cells.Item(9, 3).Value = "Ntrk3"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 3.275965666666667
$ws.Cells.Item(9, 8).Value = 9.827897
$ws.Cells.Item(9, 9).Value = 0.1819587962318212
$ws.Cells.Item(9, 10).Value = 0.1819587962318212
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.07403633333333333
$ws.Cells.Item(9, 14).Value = 0.222109
$ws.Cells.Item(9, 15).Value = 0.08034383307391152
$ws.Cells.Item(9, 16).Value = 0.08034383307391152
$ws.Cells.Item(9, 17).Value = 0.2425404860858889
$ws.Cells.Item(9, 18).Value = 2.182864374773
$ws.Cells.Item(9, 19).Value = 0.01461926715077932
$ws.Cells.Item(9, 20).Value = 0.01461926715077932

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Ntf3"
$ws.Cells.Item(10, 3).Value = "Ntrk3"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 3.275965666666667
$ws.Cells.Item(10, 8).Value = 9.827897
$ws.Cells.Item(10, 9).Value = 0.1819587962318212
$ws.Cells.Item(10, 10).Value = 0.1819587962318212
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.8384493333333333
$ws.Cells.Item(10, 14).Value = 2.515348
$ws.Cells.Item(10, 15).Value = 0.9098807334903006
$ws.Cells.Item(10, 16).Value = 0.9098807334903007
$ws.Cells.Item(10, 17).Value = 2.746731229239555
$ws.Cells.Item(10, 18).Value = 24.720581063156
$ws.Cells.Item(10, 19).Value = 0.1655608029804216
$ws.Cells.Item(10, 20).Value = 0.1655608029804216
